$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($Row1, $Row2) {
    # Columns whose values differ between the two rows and must be swapped
    $cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

    foreach ($col in $cols) {
        $addr1 = "$col$Row1"
        $addr2 = "$col$Row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }

    # Substrate info (AJ, AK, AO) moves from Row1 to Row2 (Row1 becomes blank)
    $subCols = @("AJ", "AK", "AO")
    foreach ($col in $subCols) {
        $addr1 = "$col$Row1"
        $addr2 = "$col$Row2"
        $v1 = $ws.Range($addr1).Value2
        $ws.Range($addr2).Value = $v1
        $ws.Range($addr1).Value = $null
    }
}

Swap-Rows 16 17
Swap-Rows 19 20
